# Auto-generated edit script: updates cryptos list snapshot (prices / 1h volume
# deltas) plus a same-row swap of the Toncoin / RenderToken entries (rows 29-30),
# matching the GitHub Actions scrape commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E ("Volume(1h)") updates -------------------------------------------
# These are never ambiguous with numbers (padded with spaces + "%"), so a plain
# .Value assignment keeps them as text, matching the source inlineStr cells.
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("E3").Value = '  +0.71%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E5").Value = '  +1.04%  '
$ws.Range("E6").Value = '  -1.83%  '
$ws.Range("E7").Value = '  +7.06%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +7.47%  '
$ws.Range("E10").Value = '  +7.88%  '
$ws.Range("E11").Value = '  +2.48%  '
$ws.Range("E12").Value = '  +9.81%  '
$ws.Range("E13").Value = '  +0.14%  '
$ws.Range("E14").Value = '  +1.00%  '
$ws.Range("E15").Value = '  +7.87%  '
$ws.Range("E16").Value = '  +43.34%  '
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("E18").Value = '  +5.15%  '
$ws.Range("E19").Value = '  +6.44%  '
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("E21").Value = '  +42.13%  '
$ws.Range("E22").Value = '  +9.98%  '
$ws.Range("E23").Value = '  +0.91%  '
$ws.Range("E24").Value = '  +1.87%  '
$ws.Range("E25").Value = '  +2.64%  '
$ws.Range("E26").Value = '  +11.92%  '
$ws.Range("E27").Value = '  +7.76%  '
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("E29").Value = '  -7.38%  '
$ws.Range("E30").Value = '  -1.41%  '
$ws.Range("E31").Value = '  +6.00%  '
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("E33").Value = '  +0.18%  '
$ws.Range("E34").Value = '  -3.36%  '
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("E36").Value = '  +3.52%  '
$ws.Range("E37").Value = '  +4.14%  '
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("E39").Value = '  +1.52%  '
$ws.Range("E40").Value = '  +7.77%  '
$ws.Range("E41").Value = '  -0.43%  '
$ws.Range("E42").Value = '  +1.93%  '
$ws.Range("E43").Value = '  -2.74%  '
$ws.Range("E44").Value = '  +7.42%  '
$ws.Range("E45").Value = '  +1.52%  '
$ws.Range("E46").Value = '  +13.36%  '
$ws.Range("E47").Value = '  -0.38%  '
$ws.Range("E48").Value = '  +5.93%  '
$ws.Range("E49").Value = '  +12.75%  '
$ws.Range("E50").Value = '  +0.95%  '
$ws.Range("E51").Value = '  +1.66%  '

# --- Columns B/C/D updates -------------------------------------------------------
# Several of these look like plain numbers (e.g. "42.67") or use dotted
# thousands separators ("61.817.37"), which Excel would silently coerce to a
# Number (or misparse) if assigned straight to .Value. To guarantee they land as
# text - like the original inlineStr cells - each target cell first gets a
# formula whose result is a text literal; PasteSpecial(xlPasteValues) then bakes
# that formula down to a static string in place (cell-by-cell, since a single
# multi-area Union paste does not reliably convert every area). The cell style is
# left untouched throughout (no NumberFormat side effects).
$c = $ws.Range("D2")
$c.Formula = '="61.817.37"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D3")
$c.Formula = '="3.422.22"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D5")
$c.Formula = '="408.88"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D6")
$c.Formula = '="128.24"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D8")
$c.Formula = '="0.999"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D11")
$c.Formula = '="42.67"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D12")
$c.Formula = '="9.16"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D14")
$c.Formula = '="3.964.42"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D15")
$c.Formula = '="21.32"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D16")
$c.Formula = '="0.0000204"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D17")
$c.Formula = '="3.419.63"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D18")
$c.Formula = '="12.27"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D19")
$c.Formula = '="1.07"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D20")
$c.Formula = '="61.889.36"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D21")
$c.Formula = '="441.78"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D22")
$c.Formula = '="91.44"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D24")
$c.Formula = '="12.95"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D26")
$c.Formula = '="32.97"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D27")
$c.Formula = '="8.65"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("B29")
$c.Formula = '="RenderToken"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("C29")
$c.Formula = '="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D29")
$c.Formula = '="7.60"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("B30")
$c.Formula = '="Toncoin"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("C30")
$c.Formula = '="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D30")
$c.Formula = '="2.72"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D33")
$c.Formula = '="0.115"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D34")
$c.Formula = '="42.69"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D36")
$c.Formula = '="0.0498"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D37")
$c.Formula = '="53.41"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D41")
$c.Formula = '="2.94"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D42")
$c.Formula = '="142.22"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D44")
$c.Formula = '="4.22"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D47")
$c.Formula = '="16.57"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D48")
$c.Formula = '="22.36"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D50")
$c.Formula = '="3.772.08"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$c = $ws.Range("D51")
$c.Formula = '="2.125.19"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
